$d = $word.ActiveDocument

# 1 & 6: Title / heading text appears twice (main H1 heading and the bolded
# "recap" run near the end). Neither run is adjacent to a same-formatting
# empty run that would get merged away by a plain text replace (the H1 run
# has no leading empty run; the bold run's leading empty run has different
# formatting than the text run, so it survives a normal replace), so a
# simple ReplaceAll is safe here.
$d.Content.Find.Execute("Play Faust for Free - Dark Literary-Themed Slot Game", $true, $false, $false, $false, $false, $true, 1, $false, "Play Faust Slot Game Free", 2)

# 7: Closing italic summary paragraph. Its leading empty run has different
# formatting (no rPr) than the italic text run, so it is preserved by a
# plain replace too.
$d.Content.Find.Execute("Discover the formula to defeat the devil in Faust, an immersive slot game inspired by Goethe's masterpiece. Play for free and access bonus features.", $true, $false, $false, $false, $false, $true, 1, $false, "Discover the immersive atmosphere and exciting bonus features of Faust. Play for free.", 2)

# 2-5: "What we like" / "What we don't like" bullet points. These runs have
# no rPr at all, identical to their leading empty run, so a plain text
# replace merges the two runs into one (losing the empty run). Briefly
# toggling Bold on the found range before/after the text swap forces the
# run to stay distinct from its neighboring empty run, preserving the
# paragraph's original <w:r/> + text-run shape.
function Replace-BulletText($oldText, $newText) {
    $fr = $d.Content
    $fr.Find.Execute($oldText)
    $fr.Font.Bold = 1
    $fr.Text = $newText
    $fr.Font.Bold = 0
}

Replace-BulletText "Immersive and dark atmosphere" "Immersive atmosphere with dark and eerie mood"
Replace-BulletText "Straightforward gameplay" "Simple and straightforward gameplay"
Replace-BulletText "Bonus features including free spins" "Exciting bonus features and free spins"
Replace-BulletText "No progressive jackpot" "No progressive jackpot feature"
